# Updates cryptos list data (price/volume columns, plus a few rank
# reshuffles/replacements) to match the latest coinranking.com scrape.
#
# Column D ("Price") cells are stored as text in the workbook (values like
# "26.293.79" or "1.585.74" aren't valid numbers), so for every D-column
# write we briefly force NumberFormat="@" (Text) before assigning the
# value -- otherwise Excel's COM layer auto-coerces the string into a
# number (dropping trailing zeros, e.g. "64.50" -> 64.5, or mis-parsing
# multi-dot values). Style is reset to 'Normal' right after so no extra
# cell formatting sticks around.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '26.293.79'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -1.58%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.585.74'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -0.96%  '
$ws.Range('E4').Value = '  -0.18%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '209.58'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.80%  '
$ws.Range('E6').Value = '  -1.37%  '
$ws.Range('E7').Value = '  -0.18%  '
$ws.Range('E8').Value = '  -0.97%  '
$ws.Range('E9').Value = '  -0.31%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '19.62'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.24%  '
$ws.Range('E11').Value = '  +0.26%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.809.06'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.95%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.585.59'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -1.06%  '
$ws.Range('E14').Value = '  -0.51%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.517'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -1.26%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '64.50'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -1.11%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '26.292.54'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E18').Value = '  -0.69%  '
$ws.Range('E19').Value = '  -0.04%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '207.28'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -1.76%  '
$ws.Range('E22').Value = '  -0.96%  '
$ws.Range('E23').Value = '  -3.82%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '8.84'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -1.54%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '144.36'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.54%  '
$ws.Range('E26').Value = '  -0.18%  '
$ws.Range('E27').Value = '  -1.44%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.113'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.55%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '15.28'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.81%  '
$ws.Range('E30').Value = '  -2.04%  '
$ws.Range('E31').Value = '  -0.90%  '
$ws.Range('E32').Value = '  -0.89%  '
$ws.Range('B33').Value = 'InternetComputer(DFINITY)'
$ws.Range('C33').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.95'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.87%  '
$ws.Range('B34').Value = 'WEMIXToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.31'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +12.81%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.283.59'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -1.21%  '
$ws.Range('E36').Value = '  +0.33%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.618'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +1.60%  '
$ws.Range('E38').Value = '  -1.04%  '
$ws.Range('E39').Value = '  -1.18%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.820'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.33%  '
$ws.Range('E41').Value = '  +1.02%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.770'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -1.49%  '
$ws.Range('E43').Value = '  -2.84%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '62.40'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -1.24%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.721.24'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.86%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '88.96'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -2.26%  '
$ws.Range('E47').Value = '  -0.42%  '
$ws.Range('B48').Value = 'BabyDogeCoin'
$ws.Range('C48').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0₆0103'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.99%  '
$ws.Range('B49').Value = 'Algorand'
$ws.Range('C49').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.102'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.49%  '
$ws.Range('B50').Value = 'Cronos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0510'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -1.68%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '7.48'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.59%  '
